$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Nota:" / "Fuente:" footnote rows (A10:D11) - rich-text shared strings removed
$ws.Range("A10:D11").ClearContents()

# Give the (now empty) spacer row 9 more breathing room
$ws.Rows(9).RowHeight = 21.75

# Column D is widened; A:C and E keep their original width
$ws.Columns("D").ColumnWidth = 27.5

# Move the active selection to C15
$ws.Range("C15").Select()
